$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before C: old Locator column (B) is kept as
# "LocatorName", the new column C becomes "LocatorValue", and the former
# Action/Value columns shift right from C/D to D/E.
$ws.Columns("C").Insert()

# New column matches the width of the Locator (B) column it was split from.
$ws.Columns("C").ColumnWidth() = $ws.Columns("B").ColumnWidth()

# Header row
$ws.Range("B1").Value() = "LocatorName"
$ws.Range("C1").Value() = "LocatorValue"

# Row 2: open browser - no locator
$ws.Range("B2").Value() = "NA"
$ws.Range("C2").Value() = "NA"

# Row 3: launch url - no locator
$ws.Range("B3").Value() = "NA"
$ws.Range("C3").Value() = "NA"

# Row 4: enter username - locator was "name=username"
$ws.Range("B4").Value() = "name"
$ws.Range("C4").Value() = "username"

# Row 5: enter password - locator was "name=password"
$ws.Range("B5").Value() = "name"
$ws.Range("C5").Value() = "password"

# Row 6: click login - locator was "xpath=/html/body/div/div[2]/div/div/div/div/div[2]/form/div[3]/button"
$ws.Range("B6").Value() = "xpath"
$ws.Range("C6").Value() = "//button[@type='submit']"

# Row 7: close browser - no locator (previously blank)
$ws.Range("B7").Value() = "NA"
$ws.Range("C7").Value() = "NA"

# Update the active selection to C7
$ws.Range("C7").Select()
